$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename first sheet: "cep_a_buscar" -> "cep_a_consultar"
$ws1.Name = "cep_a_consultar"

# Add the new CEPs to look up on sheet1 (rows 2-5)
$ws1.Range("A2").Value = 38401220
$ws1.Range("A3").Value = 38408240
$ws1.Range("A4").Value = 38400656
$ws1.Range("A5").Value = 38410234

# Configure sheet1 page setup (portrait, paper size 9 / A4)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Update the existing result row on "dados_coletados" with the new lookup data
$ws2.Range("A2").Value = "Rua Albertino Silva "
$ws2.Range("B2").Value = "Presidente Roosevelt "
$ws2.Range("C2").Value = "Uberlândia/MG "
$ws2.Range("D2").Value = "38401-220"

# New rows of collected data (rows 3-5), one per additional CEP
$ws2.Range("A3").Value = "Rua João Catanduva "
$ws2.Range("B3").Value = "Santa Mônica "
$ws2.Range("C3").Value = "Uberlândia/MG "
$ws2.Range("D3").Value = "38408-240"

$ws2.Range("A4").Value = "Rua São Paulo - de 157/158 a 1569/1570 "
$ws2.Range("B4").Value = "Brasil "
$ws2.Range("C4").Value = "Uberlândia/MG "
$ws2.Range("D4").Value = "38400-656"

$ws2.Range("A5").Value = "Rua Amador Lourenço "
$ws2.Range("B5").Value = "Laranjeiras "
$ws2.Range("C5").Value = "Uberlândia/MG "
$ws2.Range("D5").Value = "38410-234"

# Rows 6-9 repeat the data from rows 2-5
$ws2.Range("A6").Value = "Rua Albertino Silva "
$ws2.Range("B6").Value = "Presidente Roosevelt "
$ws2.Range("C6").Value = "Uberlândia/MG "
$ws2.Range("D6").Value = "38401-220"

$ws2.Range("A7").Value = "Rua João Catanduva "
$ws2.Range("B7").Value = "Santa Mônica "
$ws2.Range("C7").Value = "Uberlândia/MG "
$ws2.Range("D7").Value = "38408-240"

$ws2.Range("A8").Value = "Rua São Paulo - de 157/158 a 1569/1570 "
$ws2.Range("B8").Value = "Brasil "
$ws2.Range("C8").Value = "Uberlândia/MG "
$ws2.Range("D8").Value = "38400-656"

$ws2.Range("A9").Value = "Rua Amador Lourenço "
$ws2.Range("B9").Value = "Laranjeiras "
$ws2.Range("C9").Value = "Uberlândia/MG "
$ws2.Range("D9").Value = "38410-234"

# Rows 10-13 repeat the data from rows 2-5 a third time
$ws2.Range("A10").Value = "Rua Albertino Silva "
$ws2.Range("B10").Value = "Presidente Roosevelt "
$ws2.Range("C10").Value = "Uberlândia/MG "
$ws2.Range("D10").Value = "38401-220"

$ws2.Range("A11").Value = "Rua João Catanduva "
$ws2.Range("B11").Value = "Santa Mônica "
$ws2.Range("C11").Value = "Uberlândia/MG "
$ws2.Range("D11").Value = "38408-240"

$ws2.Range("A12").Value = "Rua São Paulo - de 157/158 a 1569/1570 "
$ws2.Range("B12").Value = "Brasil "
$ws2.Range("C12").Value = "Uberlândia/MG "
$ws2.Range("D12").Value = "38400-656"

$ws2.Range("A13").Value = "Rua Amador Lourenço "
$ws2.Range("B13").Value = "Laranjeiras "
$ws2.Range("C13").Value = "Uberlândia/MG "
$ws2.Range("D13").Value = "38410-234"

# Set selection on sheet2 (not the active tab after this script) then on sheet1 (active tab)
[void]$ws2.Range("E29").Select()
[void]$ws1.Range("A2").Select()

Write-Host "done"
